# Update cryptos list values (price / 1h volume change, and a few
# reordered coin rows) to match the refreshed scrape.
# Numeric-looking "Price" values are written with a leading apostrophe
# so Excel keeps them as literal text (matching the original inline
# string cells) instead of silently converting them to numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").Value = "24.438.34"
$ws.Range("E2").Value = "  +1.20%  "

$ws.Range("D3").Value = "1.667.51"
$ws.Range("E3").Value = "  +1.34%  "

$ws.Range("E4").Value = "  -0.39%  "

$ws.Range("D5").Value = "'313.57"
$ws.Range("E5").Value = "  +1.66%  "

$ws.Range("D6").Value = "'0.9998"
$ws.Range("E6").Value = "  -0.31%  "

$ws.Range("D7").Value = "'0.3966"
$ws.Range("E7").Value = "  +1.94%  "

$ws.Range("D8").Value = "'0.3919"
$ws.Range("E8").Value = "  +1.15%  "

$ws.Range("D9").Value = "'51.91"
$ws.Range("E9").Value = "  +5.81%  "

$ws.Range("D10").Value = "'1.409"
$ws.Range("E10").Value = "  +3.16%  "

$ws.Range("E11").Value = "  -0.41%  "

$ws.Range("D12").Value = "'0.08604"
$ws.Range("E12").Value = "  +1.45%  "

$ws.Range("D13").Value = "'24.50"
$ws.Range("E13").Value = "  +1.18%  "

$ws.Range("D14").Value = "'7.351"
$ws.Range("E14").Value = "  +2.56%  "

$ws.Range("B15").Value = "ShibaInu"
$ws.Range("C15").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D15").Value = "'0.00001345"
$ws.Range("E15").Value = "  +4.27%  "

$ws.Range("B16").Value = "Chainlink"
$ws.Range("C16").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D16").Value = "'7.915"
$ws.Range("E16").Value = "  +5.14%  "

$ws.Range("D17").Value = "1.665.46"
$ws.Range("E17").Value = "  +0.79%  "

$ws.Range("D18").Value = "'95.73"
$ws.Range("E18").Value = "  +0.86%  "

$ws.Range("D19").Value = "'0.06981"
$ws.Range("E19").Value = "  +0.48%  "

$ws.Range("D20").Value = "'20.62"
$ws.Range("E20").Value = "  -2.68%  "

$ws.Range("D21").Value = "'7.023"
$ws.Range("E21").Value = "  +0.68%  "

$ws.Range("D22").Value = "'0.9998"
$ws.Range("E22").Value = "  -0.26%  "

$ws.Range("E23").Value = "  -0.48%  "

$ws.Range("D24").Value = "24.427.29"
$ws.Range("E24").Value = "  +1.13%  "

$ws.Range("D25").Value = "'2.429"
$ws.Range("E25").Value = "  +3.87%  "

$ws.Range("D26").Value = "'3.036"
$ws.Range("E26").Value = "  +10.78%  "

$ws.Range("D27").Value = "'22.54"
$ws.Range("E27").Value = "  -0.19%  "

$ws.Range("D28").Value = "'157.46"
$ws.Range("E28").Value = "  -0.39%  "

$ws.Range("D29").Value = "'143.04"
$ws.Range("E29").Value = "  +0.39%  "

$ws.Range("D30").Value = "'5.400"
$ws.Range("E30").Value = "  -0.17%  "

$ws.Range("D31").Value = "'8.128"
$ws.Range("E31").Value = "  -8.34%  "

$ws.Range("D32").Value = "'2.542"
$ws.Range("E32").Value = "  +3.78%  "

$ws.Range("D33").Value = "1.855.36"
$ws.Range("E33").Value = "  +1.33%  "

$ws.Range("D34").Value = "'1.065"
$ws.Range("E34").Value = "  +7.46%  "

$ws.Range("D35").Value = "'0.08272"
$ws.Range("E35").Value = "  +2.28%  "

$ws.Range("D36").Value = "'0.03032"
$ws.Range("E36").Value = "  +2.20%  "

$ws.Range("D37").Value = "'6.825"
$ws.Range("E37").Value = "  -3.72%  "

$ws.Range("B38").Value = "FraxShare"
$ws.Range("C38").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D38").Value = "'11.12"
$ws.Range("E38").Value = "  +10.28%  "

$ws.Range("B39").Value = "Algorand"
$ws.Range("C39").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D39").Value = "'0.2766"
$ws.Range("E39").Value = "  +2.06%  "

$ws.Range("D40").Value = "'0.09267"
$ws.Range("E40").Value = "  -0.41%  "

$ws.Range("D41").Value = "'0.7753"
$ws.Range("E41").Value = "  +1.37%  "

$ws.Range("D42").Value = "'13.87"
$ws.Range("E42").Value = "  +5.73%  "

$ws.Range("D43").Value = "'1.446"
$ws.Range("E43").Value = "  -2.43%  "

$ws.Range("D44").Value = "'16.53"
$ws.Range("E44").Value = "  +1.89%  "

$ws.Range("D45").Value = "'0.7136"
$ws.Range("E45").Value = "  +3.46%  "

$ws.Range("D46").Value = "'2.541"
$ws.Range("E46").Value = "  +1.68%  "

$ws.Range("D47").Value = "'4.144"
$ws.Range("E47").Value = "  +1.25%  "

$ws.Range("D48").Value = "'0.9998"
$ws.Range("E48").Value = "  -0.29%  "

$ws.Range("D49").Value = "'0.08450"
$ws.Range("E49").Value = "  +0.30%  "

$ws.Range("B50").Value = "Tezos"
$ws.Range("C50").Value = "https://coinranking.com/coin/fsIbGOEJWbzxG+tezos-xtz"
$ws.Range("D50").Value = "'1.457"
$ws.Range("E50").Value = "  +13.96%  "

$ws.Range("B51").Value = "Quant"
$ws.Range("C51").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D51").Value = "'136.71"
$ws.Range("E51").Value = "  +1.75%  "
